# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计" (tab #2,
#    pushing every existing quarter sheet one slot to the right).
# ---------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $totalSheet)
$ws2.Name = "2022-Q4"

# Header row (B1:H1) and index column (A2:A30) reuse the same bold /
# centered / thin-border look already present on the "总计" sheet -
# copy that formatting across instead of re-deriving it by hand.
$totalSheet.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$ws2.Range("A2:A30").PasteSpecial(-4122)

# Columns B-G (rows 2-28) hold numeric-looking figures (fund codes,
# percentages, AUM, ...) that the source file stores as literal TEXT
# (t="inlineStr"), not numbers. Mark that block as Text *before*
# writing so entries like "002910" or "79.05" keep their exact
# formatting/leading zeros instead of being auto-parsed into numbers.
$ws2.Range("B2:G28").NumberFormat = "@"

$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "002910"
$ws2.Range("C2").Value = "易方达供给改革灵活配置混合"
$ws2.Range("D2").Value = "79.05"
$ws2.Range("E2").Value = "92.55"
$ws2.Range("F2").Value = "6.78"
$ws2.Range("G2").Value = "5.3596"
$ws2.Range("H2").Value = 7
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "160505"
$ws2.Range("C3").Value = "博时主题行业混合（LOF）"
$ws2.Range("D3").Value = "66.23"
$ws2.Range("E3").Value = "79.29"
$ws2.Range("F3").Value = "1.98"
$ws2.Range("G3").Value = "1.3114"
$ws2.Range("H3").Value = 8
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "398051"
$ws2.Range("C4").Value = "中海环保新能源混合"
$ws2.Range("D4").Value = "19.89"
$ws2.Range("E4").Value = "77.11"
$ws2.Range("F4").Value = "3.42"
$ws2.Range("G4").Value = "0.6802"
$ws2.Range("H4").Value = 7
$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "398021"
$ws2.Range("C5").Value = "中海能源策略混合"
$ws2.Range("D5").Value = "18.04"
$ws2.Range("E5").Value = "90.92"
$ws2.Range("F5").Value = "3.59"
$ws2.Range("G5").Value = "0.6476"
$ws2.Range("H5").Value = 5
$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "010328"
$ws2.Range("C6").Value = "博时荣华灵活配置混合A"
$ws2.Range("D6").Value = "4.01"
$ws2.Range("E6").Value = "72.54"
$ws2.Range("F6").Value = "6.31"
$ws2.Range("G6").Value = "0.2530"
$ws2.Range("H6").Value = 1
$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "005904"
$ws2.Range("C7").Value = "华泰保兴成长优选混合A"
$ws2.Range("D7").Value = "3.52"
$ws2.Range("E7").Value = "70.92"
$ws2.Range("F7").Value = "3.94"
$ws2.Range("G7").Value = "0.1387"
$ws2.Range("H7").Value = 4
$ws2.Range("A8").Value = 6
$ws2.Range("B8").Value = "003434"
$ws2.Range("C8").Value = "博时鑫泽灵活配置混合A"
$ws2.Range("D8").Value = "4.41"
$ws2.Range("E8").Value = "65.87"
$ws2.Range("F8").Value = "2.28"
$ws2.Range("G8").Value = "0.1005"
$ws2.Range("H8").Value = 6
$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = "001277"
$ws2.Range("C9").Value = "博时国企改革主题股票A"
$ws2.Range("D9").Value = "2.31"
$ws2.Range("E9").Value = "87.79"
$ws2.Range("F9").Value = "4.19"
$ws2.Range("G9").Value = "0.0968"
$ws2.Range("H9").Value = 1
$ws2.Range("A10").Value = 8
$ws2.Range("B10").Value = "050014"
$ws2.Range("C10").Value = "博时创业成长混合A"
$ws2.Range("D10").Value = "1.41"
$ws2.Range("E10").Value = "83.80"
$ws2.Range("F10").Value = "6.18"
$ws2.Range("G10").Value = "0.0871"
$ws2.Range("H10").Value = 1
$ws2.Range("A11").Value = 9
$ws2.Range("B11").Value = "014232"
$ws2.Range("C11").Value = "博时专精特新主题混合A"
$ws2.Range("D11").Value = "3.00"
$ws2.Range("E11").Value = "80.89"
$ws2.Range("F11").Value = "2.48"
$ws2.Range("G11").Value = "0.0744"
$ws2.Range("H11").Value = 1
$ws2.Range("A12").Value = 10
$ws2.Range("B12").Value = "014233"
$ws2.Range("C12").Value = "博时专精特新主题混合C"
$ws2.Range("D12").Value = "2.58"
$ws2.Range("E12").Value = "80.89"
$ws2.Range("F12").Value = "2.48"
$ws2.Range("G12").Value = "0.0640"
$ws2.Range("H12").Value = 1
$ws2.Range("A13").Value = 11
$ws2.Range("B13").Value = "002281"
$ws2.Range("C13").Value = "建信裕利灵活配置混合"
$ws2.Range("D13").Value = "0.85"
$ws2.Range("E13").Value = "93.43"
$ws2.Range("F13").Value = "5.01"
$ws2.Range("G13").Value = "0.0426"
$ws2.Range("H13").Value = 3
$ws2.Range("A14").Value = 12
$ws2.Range("B14").Value = "002378"
$ws2.Range("C14").Value = "建信弘利灵活配置混合A"
$ws2.Range("D14").Value = "0.82"
$ws2.Range("E14").Value = "92.97"
$ws2.Range("F14").Value = "4.65"
$ws2.Range("G14").Value = "0.0381"
$ws2.Range("H14").Value = 3
$ws2.Range("A15").Value = 13
$ws2.Range("B15").Value = "160519"
$ws2.Range("C15").Value = "博时睿利事件驱动灵活配置混合"
$ws2.Range("D15").Value = "0.62"
$ws2.Range("E15").Value = "87.53"
$ws2.Range("F15").Value = "5.79"
$ws2.Range("G15").Value = "0.0359"
$ws2.Range("H15").Value = 1
$ws2.Range("A16").Value = 14
$ws2.Range("B16").Value = "010690"
$ws2.Range("C16").Value = "万家互联互通核心资产量化策略混合A"
$ws2.Range("D16").Value = "0.53"
$ws2.Range("E16").Value = "88.45"
$ws2.Range("F16").Value = "6.43"
$ws2.Range("G16").Value = "0.0341"
$ws2.Range("H16").Value = 5
$ws2.Range("A17").Value = 15
$ws2.Range("B17").Value = "009967"
$ws2.Range("C17").Value = "博时荣泰灵活配置混合"
$ws2.Range("D17").Value = "0.97"
$ws2.Range("E17").Value = "82.22"
$ws2.Range("F17").Value = "3.51"
$ws2.Range("G17").Value = "0.0340"
$ws2.Range("H17").Value = 2
$ws2.Range("A18").Value = 16
$ws2.Range("B18").Value = "014999"
$ws2.Range("C18").Value = "华泰保兴吉年盈混合A"
$ws2.Range("D18").Value = "0.84"
$ws2.Range("E18").Value = "84.95"
$ws2.Range("F18").Value = "3.69"
$ws2.Range("G18").Value = "0.0310"
$ws2.Range("H18").Value = 8
$ws2.Range("A19").Value = 17
$ws2.Range("B19").Value = "015276"
$ws2.Range("C19").Value = "博时均衡回报混合A"
$ws2.Range("D19").Value = "0.91"
$ws2.Range("E19").Value = "73.67"
$ws2.Range("F19").Value = "2.28"
$ws2.Range("G19").Value = "0.0207"
$ws2.Range("H19").Value = 9
$ws2.Range("A20").Value = 18
$ws2.Range("B20").Value = "014627"
$ws2.Range("C20").Value = "财通多策略福瑞混合（LOF）C"
$ws2.Range("D20").Value = "1.09"
$ws2.Range("E20").Value = "62.26"
$ws2.Range("F20").Value = "1.20"
$ws2.Range("G20").Value = "0.0131"
$ws2.Range("H20").Value = 5
$ws2.Range("A21").Value = 19
$ws2.Range("B21").Value = "003435"
$ws2.Range("C21").Value = "博时鑫泽灵活配置混合C"
$ws2.Range("D21").Value = "0.54"
$ws2.Range("E21").Value = "65.87"
$ws2.Range("F21").Value = "2.28"
$ws2.Range("G21").Value = "0.0123"
$ws2.Range("H21").Value = 6
$ws2.Range("A22").Value = 20
$ws2.Range("B22").Value = "501028"
$ws2.Range("C22").Value = "财通多策略福瑞混合（LOF）A"
$ws2.Range("D22").Value = "0.81"
$ws2.Range("E22").Value = "62.26"
$ws2.Range("F22").Value = "1.20"
$ws2.Range("G22").Value = "0.0097"
$ws2.Range("H22").Value = 5
$ws2.Range("A23").Value = 21
$ws2.Range("B23").Value = "010691"
$ws2.Range("C23").Value = "万家互联互通核心资产量化策略混合C"
$ws2.Range("D23").Value = "0.15"
$ws2.Range("E23").Value = "88.45"
$ws2.Range("F23").Value = "6.43"
$ws2.Range("G23").Value = "0.0096"
$ws2.Range("H23").Value = 5
$ws2.Range("A24").Value = 22
$ws2.Range("B24").Value = "010329"
$ws2.Range("C24").Value = "博时荣华灵活配置混合C"
$ws2.Range("D24").Value = "0.14"
$ws2.Range("E24").Value = "72.54"
$ws2.Range("F24").Value = "6.31"
$ws2.Range("G24").Value = "0.0088"
$ws2.Range("H24").Value = 1
$ws2.Range("A25").Value = 23
$ws2.Range("B25").Value = "005905"
$ws2.Range("C25").Value = "华泰保兴成长优选混合C"
$ws2.Range("D25").Value = "0.11"
$ws2.Range("E25").Value = "70.92"
$ws2.Range("F25").Value = "3.94"
$ws2.Range("G25").Value = "0.0043"
$ws2.Range("H25").Value = 4
$ws2.Range("A26").Value = 24
$ws2.Range("B26").Value = "002553"
$ws2.Range("C26").Value = "博时创业成长混合C"
$ws2.Range("D26").Value = "0.07"
$ws2.Range("E26").Value = "83.80"
$ws2.Range("F26").Value = "6.18"
$ws2.Range("G26").Value = "0.0043"
$ws2.Range("H26").Value = 1
$ws2.Range("A27").Value = 25
$ws2.Range("B27").Value = "015277"
$ws2.Range("C27").Value = "博时均衡回报混合C"
$ws2.Range("D27").Value = "0.14"
$ws2.Range("E27").Value = "73.67"
$ws2.Range("F27").Value = "2.28"
$ws2.Range("G27").Value = "0.0032"
$ws2.Range("H27").Value = 9
$ws2.Range("A28").Value = 26
$ws2.Range("B28").Value = "015000"
$ws2.Range("C28").Value = "华泰保兴吉年盈混合C"
$ws2.Range("D28").Value = "0.02"
$ws2.Range("E28").Value = "84.95"
$ws2.Range("F28").Value = "3.69"
$ws2.Range("G28").Value = "0.0007"
$ws2.Range("H28").Value = 8
$ws2.Range("A29").Value = 27
$ws2.Range("B29").Value = "014382"
$ws2.Range("C29").Value = "博时国企改革主题股票C"
$ws2.Range("D29").Value = "0.00"
$ws2.Range("E29").Value = "87.79"
$ws2.Range("F29").Value = "4.19"
$ws2.Range("G29").Value = 0
$ws2.Range("H29").Value = 1
$ws2.Range("A30").Value = 28
$ws2.Range("B30").Value = "017194"
$ws2.Range("C30").Value = "建信弘利灵活配置混合C"
$ws2.Range("D30").Value = "0.00"
$ws2.Range("E30").Value = "92.97"
$ws2.Range("F30").Value = "4.65"
$ws2.Range("G30").Value = 0
$ws2.Range("H30").Value = 3

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    data right under the header, shifting the previously-existing
#    rows down by one (2022-Q3 -> row3, 2022-Q2 -> row4, ...,
#    2021-Q3 -> row7).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws1.Range("A$dst").Value = $ws1.Range("A$r").Value()
    $ws1.Range("B$dst").Value = $ws1.Range("B$r").Value()
    $ws1.Range("C$dst").Value = $ws1.Range("C$r").Value()
    $ws1.Range("D$dst").Value = $ws1.Range("D$r").Value()
}

# Row 7 is brand new (previously out of range) - copy the index-column
# style from a neighbouring cell so it matches A2:A6 (bold/centered/
# bordered), same as every other row in column A.
$ws1.Range("A6").Copy()
$ws1.Range("A7").PasteSpecial(-4122)

# Fix up the index column (0,1,2,3,4,5) now that everything shifted.
$ws1.Range("A2").Value = 0
$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
$ws1.Range("A6").Value = 4
$ws1.Range("A7").Value = 5

# New row 2: 2022-Q4 summary figures.
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 29
$ws1.Range("D2").Value = 9.119999999999999
